$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity corrections in the screws sub-table ---
# Row 14 (M2*10mm screw): 10 -> 15
$ws.Range("C14").Value = 15
# Row 16 (M2*3mm screw): 15 -> 10
$ws.Range("C16").Value = 10

# --- Re-apply the "highlighted" (orange) formatting used elsewhere in the
# sheet to the cells that now need it: copy formats only, leave values alone.
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Restore the selected cell ---
$ws.Range("C33").Select() | Out-Null
